# Weekly fruit/vegetable price update: a new price-record row is inserted
# at row 110 (pushing the existing rows 110-152 down to 111-153), and the
# new row is populated with this week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 110 - this shifts rows 110:152 down to
# 111:153 and carries their formatting, matching the existing style ($s="2"$
# on column D) down onto the inserted row.
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with this week's record. The
# non-numeric/static columns repeat the same values used throughout this
# price sheet (Ajo / Chino / Primera / China / Hortaliza / unit label).
$ws.Cells.Item(110, 1).Value = 4
$ws.Cells.Item(110, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(110, 3).Value = "Los Lagos"
$ws.Cells.Item(110, 4).Value = 44474
$ws.Cells.Item(110, 5).Value = 10
$ws.Cells.Item(110, 6).Value = 100112003
$ws.Cells.Item(110, 7).Value = "Ajo"
$ws.Cells.Item(110, 8).Value = "Chino"
$ws.Cells.Item(110, 9).Value = "Primera"
$ws.Cells.Item(110, 10).Value = 400
$ws.Cells.Item(110, 11).Value = 19000
$ws.Cells.Item(110, 12).Value = 20000
$ws.Cells.Item(110, 13).Value = 19500
$ws.Cells.Item(110, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(110, 15).Value = "China"
$ws.Cells.Item(110, 16).Value = 1950
$ws.Cells.Item(110, 17).Value = 10
$ws.Cells.Item(110, 18).Value = "Hortaliza"
